$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 449
$ws1.Range("F6").Value = 7661
$ws1.Range("F17").Value = 1141
$ws1.Range("F18").Value = 4596
$ws1.Range("F23").Value = 533
$ws1.Range("F24").Value = 3533
$ws1.Range("F29").Value = 3030
$ws1.Range("F31").Value = 108
$ws1.Range("F32").Value = 345
$ws1.Range("F35").Value = 486
$ws1.Range("F36").Value = 657
$ws1.Range("F40").Value = 47
$ws1.Range("F43").Value = 2919

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 38
$ws2.Range("F9").Value = 120

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1326

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1326
$ws4.Range("F6").Value = 7661
$ws4.Range("F17").Value = 1141
$ws4.Range("F18").Value = 4596
$ws4.Range("F23").Value = 533
$ws4.Range("F25").Value = 3533
$ws4.Range("F29").Value = 3030
$ws4.Range("F30").Value = 345
$ws4.Range("F34").Value = 486
$ws4.Range("F35").Value = 657
$ws4.Range("F38").Value = 38
$ws4.Range("F40").Value = 47
$ws4.Range("F43").Value = 2919
$ws4.Range("F49").Value = 120
